# "Create new test casesclear"
#
# 1. Switch workbook calculation to automatic (drops calcMode="manual").
# 2. Add two rows (a lone "placeholder" space cell) on RentMovie and MovieScreen.
# 3. Add a brand new "SettingScreen" worksheet at the end of the workbook with
#    its own Option1..Option7 / labels table.
# 4. Tidy up a couple of selections / the previously tabSelected sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Calculation mode -> Automatic ------------------------------------
$excel.Calculation = -4105   # xlCalculationAutomatic

# --- 2a. RentMovie: add a blank-looking row 9 in column B ----------------
$rentMovie = $wb.Worksheets.Item("RentMovie")
$rentMovie.Range("B9").Value = " "
$rentMovie.Range("B10").Select() | Out-Null

# --- 2b. MovieScreen: add a blank-looking row 16 in column E -------------
$movieScreen = $wb.Worksheets.Item("MovieScreen")
$movieScreen.Range("E16").Value = " "
$movieScreen.Range("A1:B1").Select() | Out-Null

# --- 3. New "SettingScreen" sheet -----------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "SettingScreen"

$ws.Range("A1").Value = "objectID"
$ws.Range("B1").Value = "Value"

$ws.Range("A2").Value = "Option1"
$ws.Range("A3").Value = "Option2"
$ws.Range("A4").Value = "Option3"
$ws.Range("A5").Value = "Option4"
$ws.Range("A6").Value = "Option5"
$ws.Range("A7").Value = "Option6"
$ws.Range("A8").Value = "Option7"

$ws.Range("B2").Value = "mijn account"
$ws.Range("B3").Value = "kinderslot"
$ws.Range("B4").Value = "aanbevelingen"
$ws.Range("B5").Value = "systeem"
$ws.Range("B6").Value = "zenders hernummeren"
$ws.Range("B7").Value = "mijn voorkeur"
$ws.Range("B8").Value = "beheer van de toestellen"

# Copy cell formatting from existing cells so the new sheet reuses the same
# shared style entries instead of creating duplicates.
$styleHeader = $wb.Worksheets.Item("RentMovie").Range("A1")
$styleHeader.Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null

$styleLabel = $wb.Worksheets.Item("MiniEPGScreen").Range("A6")
$styleLabel.Copy() | Out-Null
$ws.Range("A3:A8").PasteSpecial(-4122) | Out-Null

$styleValue = $wb.Worksheets.Item("SystemInfoScreen").Range("B2")
$styleValue.Copy() | Out-Null
$ws.Range("B2:B8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("H8").Select() | Out-Null
